# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 792
$ws1.Range("F5").Value = 151
$ws1.Range("F8").Value = 359
$ws1.Range("F9").Value = 471
$ws1.Range("F10").Value = 518
$ws1.Range("F12").Value = 12017
$ws1.Range("F13").Value = 5440

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 792
$ws4.Range("F7").Value = 151
$ws4.Range("F10").Value = 359
$ws4.Range("F11").Value = 471
$ws4.Range("F12").Value = 518
$ws4.Range("F14").Value = 12017
$ws4.Range("F16").Value = 5440

$wb.Save()
